$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.867.39"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.599.43"
$ws.Range("E3").Value = "  -2.04%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.11"
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.478"
$ws.Range("E7").Value = "  -5.06%  "
$ws.Range("E8").Value = "  -2.63%  "
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.81"
$ws.Range("E10").Value = "  -3.60%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "1.821.35"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").Value = "1.592.20"
$ws.Range("E13").Value = "  -4.24%  "
$ws.Range("E14").Value = "  -3.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.507"
$ws.Range("E15").Value = "  -4.37%  "
$ws.Range("D16").Value = "25.859.85"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.51"
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").Value = "0.0₃0706"
$ws.Range("E18").Value = "  -5.18%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "189.26"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.16"
$ws.Range("E21").Value = "  -1.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.28"
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.92"
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.71"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("E26").Value = "  -3.87%  "
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("E28").Value = "  -3.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.91"
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0462"
$ws.Range("E31").Value = "  -4.29%  "
$ws.Range("E32").Value = "  -2.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  -4.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.39"
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("D36").Value = "1.103.23"
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.792"
$ws.Range("E38").Value = "  -8.48%  "
$ws.Range("E39").Value = "  -2.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.493"
$ws.Range("E40").Value = "  -5.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "95.57"
$ws.Range("E41").Value = "  -3.06%  "
$ws.Range("D42").Value = "1.734.17"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("E43").Value = "  -3.69%  "
$ws.Range("E44").Value = "  -5.10%  "
$ws.Range("D45").Value = "0.0₆0107"
$ws.Range("E45").Value = "  -5.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "52.96"
$ws.Range("E46").Value = "  -3.77%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0509"
$ws.Range("E47").Value = "  -3.18%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.44"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("E49").Value = "  -0.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.01"
$ws.Range("E50").Value = "  -0.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.36"
$ws.Range("E51").Value = "  -1.85%  "
